# CORRIDAS TK MULTIMARCAS.xlsx - add a missing "RAIZAL" trip entry on sheet "MES 01"
# A new data row (value 10, destination "RAIZAL") is inserted right above the
# "45216" date-separator row (old row 57), pushing that row and everything
# below it down by one. The TOTAL formula at the bottom is then repaired so it
# sums the new row instead of the old broken #REF! placeholder, and the
# now-obsolete underline on the "SAO MIGUEL" entry above it is cleared.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MES 01")
$ws.Activate()

# Insert a new row above the old row 57 (the "45216" merged date cell),
# shifting it and all following rows down by one.
$ws.Rows.Item(57).Insert()

# Copy the formatting of the row above (row 54: a normal "value / place" row)
# into the freshly-inserted row 57 so it matches the other entries.
$ws.Range("A54:B54").Copy()
$ws.Range("A57:B57").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new trip: 10 minutes to RAIZAL.
$ws.Range("A57").Value = 10
$ws.Range("B57").Value = "RAIZAL"

# Repair the TOTAL formula (now on row 68): replace the dangling #REF! with
# the newly populated A56:A57 range, and shift the remaining day totals down
# one row to match the inserted row.
$ws.Range("B68").Formula = "=SUM(A3,A5,A7,A9:A13,A15,A17:A20,A22:A25,A27:A28,A30:A31,A33:A39,A41:A43,A45:A47,A49:A52,A54,A56:A57,A59,A61,A63,A65,A67)"

# The "SAO MIGUEL" entry (row 56) no longer needs its underline.
$ws.Range("B56").Font.Underline = $false

# Restore the on-screen selection/scroll position left by the edit.
$ws.Range("I60").Select()
